$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update source task hours (Day 4 column = F) for several tasks, and the
# "Basement Room A: Event" estimated hours (B20), per the 20-03 meeting update.
$ws.Range("F7").Value = 1
$ws.Range("F12").Value = 3
$ws.Range("F14").Value = 2
$ws.Range("B20").Value = 1
$ws.Range("F20").Value = 0.16

# Update the active selection to reflect where the user left off.
$ws.Range("F10").Select()
